$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2479.4517
$ws.Range("I98").Value = 2656.5417
$ws.Range("J98").Value = 1872.2858
$ws.Range("K98").Value = 2656.5417
$ws.Range("L98").Value = 1872.2858
$ws.Range("M98").Value = -1158.5417
$ws.Range("N98").Value = -4868.2858

$ws.Range("H122").Value = 2479.4517
$ws.Range("I122").Value = 2656.5417
$ws.Range("J122").Value = 1872.2858
$ws.Range("K122").Value = 7969.625100000001
$ws.Range("L122").Value = 5616.857400000001
$ws.Range("M122").Value = -5519.625100000001
$ws.Range("N122").Value = -10516.8574

$ws.Range("H132").Value = 273729.3
$ws.Range("I132").Value = 306778.3
$ws.Range("J132").Value = 1075
$ws.Range("K132").Value = 920334.8999999999
$ws.Range("L132").Value = 3225
$ws.Range("M132").Value = -917804.8999999999
$ws.Range("N132").Value = -8285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4753.254
$ws.Range("I32").Value = 5124.654
$ws.Range("J32").Value = 3465.7334
$ws.Range("K32").Value = 5124.654
$ws.Range("L32").Value = 3465.7334
$ws.Range("M32").Value = -4837.654
$ws.Range("N32").Value = -4039.7334

$ws.Range("H122").Value = 1825.8462
$ws.Range("I122").Value = 1866.909
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 5600.727000000001
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -3150.727000000001
$ws.Range("N122").Value = -9700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 733.8125
$ws.Range("I80").Value = 259.85715
$ws.Range("J80").Value = 1102.4445
$ws.Range("K80").Value = 259.85715
$ws.Range("L80").Value = 1102.4445
$ws.Range("M80").Value = 738.14285
$ws.Range("N80").Value = -3098.4445

$ws.Range("H83").Value = 733.8125
$ws.Range("I83").Value = 259.85715
$ws.Range("J83").Value = 1102.4445
$ws.Range("K83").Value = 1299.28575
$ws.Range("L83").Value = 5512.2225
$ws.Range("M83").Value = 3692.71425
$ws.Range("N83").Value = -15496.2225

$ws.Range("H107").Value = 232187
$ws.Range("I107").Value = 328771.6
$ws.Range("J107").Value = 1457.1666
$ws.Range("K107").Value = 328771.6
$ws.Range("L107").Value = 1457.1666
$ws.Range("M107").Value = -326851.6
$ws.Range("N107").Value = -5297.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1612.75
$ws.Range("I31").Value = 1222.3077
$ws.Range("J31").Value = 3304.6667
$ws.Range("K31").Value = 1222.3077
$ws.Range("L31").Value = 3304.6667
$ws.Range("M31").Value = -927.3077000000001
$ws.Range("N31").Value = -3894.6667

$ws.Range("H34").Value = 1612.75
$ws.Range("I34").Value = 1222.3077
$ws.Range("J34").Value = 3304.6667
$ws.Range("K34").Value = 1222.3077
$ws.Range("L34").Value = 3304.6667
$ws.Range("M34").Value = -1020.3077
$ws.Range("N34").Value = -3708.6667

$ws.Range("H94").Value = 2128
$ws.Range("I94").Value = 2512
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 2512
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -2061
$ws.Range("N94").Value = -2902

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 14706795
$ws.Range("I5").Value = 920.72
$ws.Range("K5").Value = 2762.16
$ws.Range("M5").Value = -2650.16

$ws.Range("H57").Value = 2000
$ws.Range("J57").Value = 3000
$ws.Range("L57").Value = 9000
$ws.Range("N57").Value = -10118

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = ""
$ws.Range("N74").Value = 0

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = ""
$ws.Range("N77").Value = 0

$ws.Range("H98").Value = 1010.5714
$ws.Range("I98").Value = 1306.6666
$ws.Range("J98").Value = 788.5
$ws.Range("K98").Value = 3919.9998
$ws.Range("L98").Value = 2365.5
$ws.Range("M98").Value = -2421.9998
$ws.Range("N98").Value = -5361.5

$ws.Range("H99").Value = 2557.1428
$ws.Range("I99").Value = 1300.5
$ws.Range("J99").Value = 4232.6665
$ws.Range("K99").Value = 3901.5
$ws.Range("L99").Value = 12697.9995
$ws.Range("M99").Value = -1655.5
$ws.Range("N99").Value = -17189.9995

$ws.Range("H113").Value = 540.8372000000001
$ws.Range("I113").Value = 520.8570999999999
$ws.Range("J113").Value = 550.4828
$ws.Range("K113").Value = 1562.5713
$ws.Range("L113").Value = 1651.4484
$ws.Range("M113").Value = 607.4287000000002
$ws.Range("N113").Value = -5991.4484

$ws.Range("H133").Value = 6097.0645
$ws.Range("I133").Value = 3042
$ws.Range("K133").Value = 9126
$ws.Range("M133").Value = -4066

$ws.Range("H135").Value = 14706795
$ws.Range("I135").Value = 920.72
$ws.Range("K135").Value = 8286.48
$ws.Range("M135").Value = -5751.48

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 20000
$ws.Range("J93").Value = 20000
$ws.Range("L93").Value = 20000
$ws.Range("N93").Value = -23744

$ws.Range("H102").Value = 1162.2069
$ws.Range("I102").Value = 949.5454999999999
$ws.Range("J102").Value = 1830.5714
$ws.Range("K102").Value = 949.5454999999999
$ws.Range("L102").Value = 1830.5714
$ws.Range("M102").Value = 672.4545000000001
$ws.Range("N102").Value = -5074.5714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1591.9574
$ws.Range("I7").Value = 1466.4
$ws.Range("J7").Value = 1958.1666
$ws.Range("K7").Value = 1466.4
$ws.Range("L7").Value = 1958.1666
$ws.Range("M7").Value = -1354.4
$ws.Range("N7").Value = -2182.1666

$ws.Range("H40").Value = 2388.5557
$ws.Range("I40").Value = 2312.125
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2312.125
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2176.125
$ws.Range("N40").Value = -3272

$ws.Range("H82").Value = 1520.0714
$ws.Range("I82").Value = 1280.091
$ws.Range("J82").Value = 2400
$ws.Range("K82").Value = 1280.091
$ws.Range("L82").Value = 2400
$ws.Range("M82").Value = -919.0909999999999
$ws.Range("N82").Value = -3122

$ws.Range("H85").Value = 1520.0714
$ws.Range("I85").Value = 1280.091
$ws.Range("J85").Value = 2400
$ws.Range("K85").Value = 1280.091
$ws.Range("L85").Value = 2400
$ws.Range("M85").Value = -32.09099999999989
$ws.Range("N85").Value = -4896

$ws.Range("H122").Value = 5808.0884
$ws.Range("I122").Value = 6469.5186
$ws.Range("J122").Value = 3256.8572
$ws.Range("K122").Value = 19408.5558
$ws.Range("L122").Value = 9770.571599999999
$ws.Range("M122").Value = -16958.5558
$ws.Range("N122").Value = -14670.5716

$ws.Range("H126").Value = 1591.9574
$ws.Range("I126").Value = 1466.4
$ws.Range("J126").Value = 1958.1666
$ws.Range("K126").Value = 4399.200000000001
$ws.Range("L126").Value = 5874.4998
$ws.Range("M126").Value = -1929.200000000001
$ws.Range("N126").Value = -10814.4998

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = ""
$ws.Range("N129").Value = 0

$ws.Range("H132").Value = 126955.78
$ws.Range("I132").Value = 225120.8
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 675362.3999999999
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -672832.3999999999
$ws.Range("N132").Value = -17808.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 230136.36
$ws.Range("J18").Value = 230136.36
$ws.Range("L18").Value = 230136.36
$ws.Range("N18").Value = -230482.36
